# Version 0.1.1.9 - The price tag format editor has been fully implemented,
# and the interface has been improved, making it more modern and concise.
#
# Data change: combined "gender + size" text values (e.g. "муж 9", "жен 7,2")
# are split into a separate gender value in column G (Пол) and a numeric
# size value in column K (Размер). Also a couple of quantities (column O)
# were corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "жен 7,2" -> gender "жен" (size was already numeric 9);
# quantity (O4) corrected 1 -> 11
$ws.Range("G4").Value = "жен"
$ws.Range("O4").Value = 11

# Row 5: "муж 9" -> gender "муж", size -> 3
$ws.Range("G5").Value = "муж"
$ws.Range("K5").Value = 3

# Row 7: "муж 8" -> gender "муж", size -> 7.5; quantity (O7) corrected 2 -> 23
$ws.Range("G7").Value = "муж"
$ws.Range("K7").Value = 7.5
$ws.Range("O7").Value = 23

# Row 9: "муж 22" -> gender "муж", size -> 2
$ws.Range("G9").Value = "муж"
$ws.Range("K9").Value = 2

# Row 10: size -> 18 (gender was already "муж")
$ws.Range("K10").Value = 18

# Reflect the updated selection in the sheet view (active cell moved to M4)
$ws.Range("M4").Select()
